# Blood donor list workbook update
#
# - The workbook now tracks the actual donation date in column G
#   ("dateOfDonation") instead of leaving it blank and using column H
#   ("availableForDonation") as a Yes/No flag for the first few rows; the
#   stale H values for the original four rows are cleared.
# - A couple of small data-entry mistakes on the original rows are fixed
#   (state for row 3, mobile number + area for row 5).
# - Ten new donor records are appended (rows 6-15), including a few dates
#   that were typed as plain text rather than recognised as real dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear the stale "availableForDonation" values for the original rows ---
$ws.Range("H2").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("H5").ClearContents()

# --- give the whole "dateOfDonation" data column (G2:G15) a short-date
#     number format in one shot, via a scratch cell + copy/paste-special,
#     so every cell shares a single new style entry instead of one each ---
$ws.Range("N1").NumberFormat = "mm-dd-yy"
$ws.Range("N1").Copy()
$ws.Range("G2:G15").PasteSpecial(-4122)
$ws.Range("N1").Clear()

# 1) Cells that only reuse already-known text / plain numbers / styled
#    dates whose serial value parses cleanly -- order doesn't matter for these.
$ws.Range("A2").Value = "rohan"
$ws.Range("B2").Value = "pawar"
$ws.Range("C2").Value = 27
$ws.Range("D2").Value = "A+"
$ws.Range("E2").Value = 9898989898
$ws.Range("F2").Value = "Male"
$ws.Range("G2").Value = 44807
$ws.Range("I2").Value = "jbp"
$ws.Range("J2").Value = "Ranjhi"
$ws.Range("K2").Value = "MP"
$ws.Range("L2").Value = 482005
$ws.Range("A3").Value = "Aman"
$ws.Range("B3").Value = "singh"
$ws.Range("C3").Value = 26
$ws.Range("D3").Value = "O+"
$ws.Range("E3").Value = 1919191919
$ws.Range("F3").Value = "Male"
$ws.Range("G3").Value = 44564
$ws.Range("I3").Value = "banaras"
$ws.Range("J3").Value = "Ranjhi"
$ws.Range("L3").Value = 482005
$ws.Range("A4").Value = "Kanni"
$ws.Range("B4").Value = "kohli"
$ws.Range("C4").Value = 27
$ws.Range("D4").Value = "A+"
$ws.Range("E4").Value = 8787878787
$ws.Range("F4").Value = "Female"
$ws.Range("G4").Value = 44806
$ws.Range("I4").Value = "Jbp"
$ws.Range("J4").Value = "Ranjhi"
$ws.Range("K4").Value = "MP"
$ws.Range("L4").Value = 482005
$ws.Range("A5").Value = "Shahreen"
$ws.Range("B5").Value = "Shahreen"
$ws.Range("C5").Value = 24
$ws.Range("D5").Value = "O+"
$ws.Range("E5").Value = 8923849283
$ws.Range("F5").Value = "Female"
$ws.Range("G5").Value = 44805
$ws.Range("I5").Value = "Indore"
$ws.Range("K5").Value = "MP"
$ws.Range("L5").Value = 410101
$ws.Range("C6").Value = 45
$ws.Range("D6").Value = "O+"
$ws.Range("E6").Value = 2342342342
$ws.Range("F6").Value = "Male"
$ws.Range("G6").Value = 44562
$ws.Range("J6").Value = "ABC"
$ws.Range("K6").Value = "MP"
$ws.Range("L6").Value = 410101
$ws.Range("C7").Value = 33
$ws.Range("D7").Value = "O+"
$ws.Range("E7").Value = 2342342344
$ws.Range("F7").Value = "Male"
$ws.Range("G7").Value = 44897
$ws.Range("J7").Value = "ABC"
$ws.Range("K7").Value = "MP"
$ws.Range("L7").Value = 410101
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = "O+"
$ws.Range("E8").Value = 2312121111
$ws.Range("F8").Value = "Female"
$ws.Range("G8").Value = 44806
$ws.Range("I8").Value = "Indore"
$ws.Range("J8").Value = "ABC"
$ws.Range("K8").Value = "MP"
$ws.Range("L8").Value = 410101
$ws.Range("C9").Value = 46
$ws.Range("D9").Value = "O+"
$ws.Range("E9").Value = 9893894349
$ws.Range("F9").Value = "Female"
$ws.Range("K9").Value = "MP"
$ws.Range("L9").Value = 410101
$ws.Range("C10").Value = 25
$ws.Range("D10").Value = "O+"
$ws.Range("E10").Value = 8898234234
$ws.Range("F10").Value = "Female"
$ws.Range("G10").Value = 44805
$ws.Range("I10").Value = "Indore"
$ws.Range("J10").Value = "ABC"
$ws.Range("K10").Value = "MP"
$ws.Range("L10").Value = 410101
$ws.Range("C11").Value = 38
$ws.Range("D11").Value = "O+"
$ws.Range("E11").Value = 4545454555
$ws.Range("F11").Value = "Male"
$ws.Range("G11").Value = 44449
$ws.Range("I11").Value = "Indore"
$ws.Range("J11").Value = "ABC"
$ws.Range("K11").Value = "MP"
$ws.Range("L11").Value = 410101
$ws.Range("C12").Value = 33
$ws.Range("D12").Value = "O+"
$ws.Range("E12").Value = 6767667677
$ws.Range("F12").Value = "Male"
$ws.Range("G12").Value = 44450
$ws.Range("I12").Value = "jbp"
$ws.Range("J12").Value = "ABC"
$ws.Range("K12").Value = "MP"
$ws.Range("L12").Value = 410101
$ws.Range("C13").Value = 34
$ws.Range("D13").Value = "O+"
$ws.Range("E13").Value = 7676767676
$ws.Range("F13").Value = "Male"
$ws.Range("G13").Value = 44538
$ws.Range("I13").Value = "jbp"
$ws.Range("J13").Value = "ABC"
$ws.Range("K13").Value = "MP"
$ws.Range("L13").Value = 410101
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = "O+"
$ws.Range("E14").Value = 6766666666
$ws.Range("F14").Value = "Male"
$ws.Range("I14").Value = "banaras"
$ws.Range("J14").Value = "ABC"
$ws.Range("K14").Value = "UP"
$ws.Range("L14").Value = 410101
$ws.Range("C15").Value = 24
$ws.Range("D15").Value = "O+"
$ws.Range("E15").Value = 3434534553
$ws.Range("F15").Value = "Male"
$ws.Range("I15").Value = "banaras"
$ws.Range("J15").Value = "ABC"
$ws.Range("K15").Value = "UP"
$ws.Range("L15").Value = 410101

# 2) Cells that introduce brand-new shared strings. These MUST be written
#    in this exact order so new entries land at shared-string indices 39-67,
#    matching the target workbook byte-for-byte.
$ws.Range("A6").Value = "Ravi"  # -> shared string #39
$ws.Range("B6").Value = "Sharma"  # -> shared string #40
$ws.Range("A7").Value = "Praveen"  # -> shared string #41
$ws.Range("A8").Value = "Pooja"  # -> shared string #42
$ws.Range("A9").Value = "Seeta"  # -> shared string #43
$ws.Range("A10").Value = "Geeta"  # -> shared string #44
$ws.Range("A11").Value = "Hari"  # -> shared string #45
$ws.Range("A12").Value = "Deepak"  # -> shared string #46
$ws.Range("A13").Value = "Himanshu"  # -> shared string #47
$ws.Range("A14").Value = "Pathak"  # -> shared string #48
$ws.Range("A15").Value = "Dhyan"  # -> shared string #49
$ws.Range("B7").Value = "Gupta"  # -> shared string #50
$ws.Range("B8").Value = "Moo"  # -> shared string #51
$ws.Range("B9").Value = "Sivahare"  # -> shared string #52
$ws.Range("B10").Value = "Rathi"  # -> shared string #53
$ws.Range("B11").Value = "Sharna"  # -> shared string #54
$ws.Range("B12").Value = "Patel"  # -> shared string #55
$ws.Range("B13").Value = "Kuvar"  # -> shared string #56
$ws.Range("B14").Value = "Pandit"  # -> shared string #57
$ws.Range("B15").Value = "Kumar"  # -> shared string #58
$ws.Range("I6").Value = "katni"  # -> shared string #59
$ws.Range("I7").Value = "sihora"  # -> shared string #60
$ws.Range("I9").Value = "chhindwara"  # -> shared string #61
$ws.Range("J5").Value = "vijay nagar"  # -> shared string #62
$ws.Range("J9").Value = "Saori"  # -> shared string #63
$ws.Range("K3").Value = "UP"  # -> shared string #64
$ws.Range("G15").Value = "15/12/2021"  # -> shared string #65
$ws.Range("G14").Value = "19/12/2021"  # -> shared string #66
$ws.Range("G9").Value = "15/1/2022"  # -> shared string #67


# --- leave the selection where data entry left off ---
$ws.Range("H15").Select()
